$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Re-order / update the vendor list in column A, and append "Microsoft" as a new row.
$values = @("Autodesk", "Amazon", "Samsung", "Microsoft Office", "Qualcomm", "Google", "Adobe", "Cisco", "IBM", "Microsoft", "Oracle")

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Autofit column A to match the widened "Microsoft Office" content.
$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null

# Update the selection to match the author's last-active cell.
$ws.Range("H10").Select() | Out-Null

$wb.Save()
